$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hindi")
$ws.Activate()

$ws.Range("B3").Value = @'
कोविड 19 के सबसे ज्यादा सामान्य लक्षण बुखार, थकान तथा सूखी खांसी हैं । कुछ रोगियों को शरीर में दर्द, नाक बंद होना, नाक बहना यानि जुकाम, गले में दर्द या दस्त हो सकता है । आरंभ में ये लक्षण ज्यादातर हल्के रूप में सामने आते हैं और धीरे-धीरे बढ़ने लगते हैं । कुछ लोग संक्रमित होते हैं, फिरभी उनमें कोई लक्षण दिखाई नहीं देता है - लेकिन वे स्वयं को अच्छा महसूस भी नहीं करते हैं । अधिकांश लोग (करीब 80%) बिना विशेष उपचार के ही इस रोग से ठीक हो जाते हैं । कोविड 19 से संक्रमित 6 में से लगभग 1व्यक्ति ही गंभीर रूप से बीमार होते हैं और उन्हें सांस लेने में कठिनाई होती है । वृद्ध व्यक्ति तथा वे लोग ही ज्यादा गंभीर रूप से बीमार होते हैं, जिन्हें हाई ब्लड-प्रेशर, हृदय संबंधी बीमारियां या डायबिटीज है । बुखार, कफ तथा सांस लेने में कष्ट होने पर लोगों को मेडिकल उपचार कराना चाहिए ।
'@

$ws.Range("B10").Select()
